$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2").Value = "moh6q"
$ws.Range("B2").Value = "moh79@qh.com"
$ws.Range("A3").Value = "moh7q"
$ws.Range("B3").Value = "hen97@qh.com"

$ws.Range("B9").Select()
